$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the outlier flag columns, matching the style of the
# existing header row (bold, centered, bordered) by copying E1's format.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Default all outlier flags to FALSE for data rows 2-23.
$ws.Range("F2:H23").Value = $false

# Row 21 (Hb 76) is flagged as a KNN outlier.
$ws.Range("F21").Value = $true
